# Auto-generated edit script applying numeric updates to the Leve profit tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2762.7778
$ws.Range("I19").Value = 2795
$ws.Range("K19").Value = 2795
$ws.Range("M19").Value = -2620
$ws.Range("H28").Value = 2885.1
$ws.Range("I28").Value = 1756.6
$ws.Range("J28").Value = 4013.6
$ws.Range("K28").Value = 1756.6
$ws.Range("L28").Value = 4013.6
$ws.Range("M28").Value = -1271.6
$ws.Range("N28").Value = -4983.6
$ws.Range("H34").Value = 5192.4546
$ws.Range("I34").Value = 5192.4546
$ws.Range("K34").Value = 5192.4546
$ws.Range("M34").Value = -4989.4546
$ws.Range("H36").Value = 5192.4546
$ws.Range("I36").Value = 5192.4546
$ws.Range("K36").Value = 5192.4546
$ws.Range("M36").Value = -4477.4546
$ws.Range("H40").Value = 5419.778
$ws.Range("I40").Value = 5955.8
$ws.Range("K40").Value = 5955.8
$ws.Range("M40").Value = -5780.8
$ws.Range("H70").Value = 6443.4443
$ws.Range("I70").Value = 999
$ws.Range("J70").Value = 10799
$ws.Range("K70").Value = 2997
$ws.Range("L70").Value = 32397
$ws.Range("M70").Value = -2727
$ws.Range("N70").Value = -32937
$ws.Range("H73").Value = 6443.4443
$ws.Range("I73").Value = 999
$ws.Range("J73").Value = 10799
$ws.Range("K73").Value = 2997
$ws.Range("L73").Value = 32397
$ws.Range("M73").Value = -2061
$ws.Range("N73").Value = -34269
$ws.Range("H116").Value = 5032.5835
$ws.Range("I116").Value = 5137.375
$ws.Range("J116").Value = 4823
$ws.Range("K116").Value = 5137.375
$ws.Range("L116").Value = 4823
$ws.Range("M116").Value = -1695.375
$ws.Range("N116").Value = -11707
$ws.Range("H135").Value = 1391.3462
$ws.Range("I135").Value = 948.5217
$ws.Range("K135").Value = 8536.695299999999
$ws.Range("M135").Value = -6001.695299999999
$ws.Range("H136").Value = 89933
$ws.Range("J136").Value = 89933
$ws.Range("L136").Value = 89933
$ws.Range("N136").Value = -100133
$ws.Range("H137").Value = 2220.9333
$ws.Range("I137").Value = 2182
$ws.Range("J137").Value = 2246.889
$ws.Range("K137").Value = 6546
$ws.Range("L137").Value = 6740.667
$ws.Range("M137").Value = -3996
$ws.Range("N137").Value = -11840.667
$ws.Range("H138").Value = 3265.4792
$ws.Range("J138").Value = 9710.615
$ws.Range("L138").Value = 29131.845
$ws.Range("N138").Value = -39411.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2575.8125
$ws.Range("I2").Value = 2851.6428
$ws.Range("K2").Value = 2851.6428
$ws.Range("M2").Value = -2738.6428
$ws.Range("H116").Value = 2575.8125
$ws.Range("I116").Value = 2851.6428
$ws.Range("K116").Value = 2851.6428
$ws.Range("M116").Value = -557.6428000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2575.8125
$ws.Range("I3").Value = 2851.6428
$ws.Range("K3").Value = 2851.6428
$ws.Range("M3").Value = -2737.6428
$ws.Range("H20").Value = 3124.5
$ws.Range("I20").Value = 2964.8333
$ws.Range("J20").Value = 3843
$ws.Range("K20").Value = 2964.8333
$ws.Range("L20").Value = 3843
$ws.Range("M20").Value = -2717.8333
$ws.Range("N20").Value = -4337
$ws.Range("H137").Value = 79764.336
$ws.Range("J137").Value = 79764.336
$ws.Range("L137").Value = 79764.336
$ws.Range("N137").Value = -89964.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8825.27
$ws.Range("I31").Value = 9387.950000000001
$ws.Range("J31").Value = 6949.6665
$ws.Range("K31").Value = 9387.950000000001
$ws.Range("L31").Value = 6949.6665
$ws.Range("M31").Value = -9092.950000000001
$ws.Range("N31").Value = -7539.6665
$ws.Range("H34").Value = 8825.27
$ws.Range("I34").Value = 9387.950000000001
$ws.Range("J34").Value = 6949.6665
$ws.Range("K34").Value = 9387.950000000001
$ws.Range("L34").Value = 6949.6665
$ws.Range("M34").Value = -9185.950000000001
$ws.Range("N34").Value = -7353.6665
$ws.Range("H39").Value = 6770
$ws.Range("I39").Value = 6770
$ws.Range("K39").Value = 6770
$ws.Range("M39").Value = -6379
$ws.Range("H49").Value = 6770
$ws.Range("I49").Value = 6770
$ws.Range("K49").Value = 6770
$ws.Range("M49").Value = -6588
$ws.Range("H135").Value = 39999
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9780.166999999999
$ws.Range("I56").Value = 9780.166999999999
$ws.Range("K56").Value = 9780.166999999999
$ws.Range("M56").Value = -9250.166999999999
$ws.Range("H107").Value = 1311.091
$ws.Range("J107").Value = 623
$ws.Range("L107").Value = 1869
$ws.Range("N107").Value = -5709
$ws.Range("H132").Value = 1941.25
$ws.Range("I132").Value = 1974.75
$ws.Range("K132").Value = 17772.75
$ws.Range("M132").Value = -15242.75
$ws.Range("H139").Value = 6973.077
$ws.Range("I139").Value = 5706.625
$ws.Range("K139").Value = 17119.875
$ws.Range("M139").Value = -11979.875
$ws.Range("H141").Value = 3999
$ws.Range("I141").Value = 3998
$ws.Range("K141").Value = 11994
$ws.Range("M141").Value = -6814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10543463
$ws.Range("I11").Value = 7583735.5
$ws.Range("K11").Value = 7583735.5
$ws.Range("M11").Value = -7583596.5
$ws.Range("H14").Value = 1839218.2
$ws.Range("I14").Value = 3200733.8
$ws.Range("K14").Value = 3200733.8
$ws.Range("M14").Value = -3200565.8
$ws.Range("H34").Value = 32079
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 32079
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 32079
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -32615
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H76").Value = 32079
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 32079
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 32079
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -32709
$ws.Range("H79").Value = 32079
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 32079
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 32079
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -34263
$ws.Range("H122").Value = 2087.6316
$ws.Range("I122").Value = 2086.1765
$ws.Range("K122").Value = 6258.529500000001
$ws.Range("M122").Value = -3808.529500000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H139").Value = 66971.25
$ws.Range("J139").Value = 66971.25
$ws.Range("L139").Value = 66971.25
$ws.Range("N139").Value = -77251.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1903.25
$ws.Range("I22").Value = 1817.3529
$ws.Range("J22").Value = 2000.6
$ws.Range("K22").Value = 1817.3529
$ws.Range("L22").Value = 2000.6
$ws.Range("M22").Value = -1522.3529
$ws.Range("N22").Value = -2590.6
$ws.Range("H27").Value = 1903.25
$ws.Range("I27").Value = 1817.3529
$ws.Range("J27").Value = 2000.6
$ws.Range("K27").Value = 1817.3529
$ws.Range("L27").Value = 2000.6
$ws.Range("M27").Value = -1710.3529
$ws.Range("N27").Value = -2214.6
$ws.Range("H82").Value = 948.2857
$ws.Range("I82").Value = 728.1429000000001
$ws.Range("J82").Value = 1168.4286
$ws.Range("K82").Value = 728.1429000000001
$ws.Range("L82").Value = 1168.4286
$ws.Range("M82").Value = -367.1429000000001
$ws.Range("N82").Value = -1890.4286
$ws.Range("H85").Value = 948.2857
$ws.Range("I85").Value = 728.1429000000001
$ws.Range("J85").Value = 1168.4286
$ws.Range("K85").Value = 728.1429000000001
$ws.Range("L85").Value = 1168.4286
$ws.Range("M85").Value = 519.8570999999999
$ws.Range("N85").Value = -3664.4286
$ws.Range("H109").Value = 532000
$ws.Range("J109").Value = 532000
$ws.Range("L109").Value = 532000
$ws.Range("N109").Value = -534774
$ws.Range("H134").Value = 85966.336
$ws.Range("J134").Value = 85966.336
$ws.Range("L134").Value = 85966.336
$ws.Range("N134").Value = -96106.336
$ws.Range("H135").Value = 86932.11
$ws.Range("J135").Value = 89311.625
$ws.Range("L135").Value = 89311.625
$ws.Range("N135").Value = -99451.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60537.832
$ws.Range("I46").Value = 76000
$ws.Range("J46").Value = 57445.4
$ws.Range("K46").Value = 76000
$ws.Range("L46").Value = 57445.4
$ws.Range("M46").Value = -75769
$ws.Range("N46").Value = -57907.4
$ws.Range("H57").Value = 103899.664
$ws.Range("J57").Value = 103899.664
$ws.Range("L57").Value = 103899.664
$ws.Range("N57").Value = -105407.664
$ws.Range("H130").Value = 104500
$ws.Range("I130").Value = 104000
$ws.Range("K130").Value = 104000
$ws.Range("M130").Value = -98980
$ws.Range("H134").Value = 60537.832
$ws.Range("I134").Value = 76000
$ws.Range("J134").Value = 57445.4
$ws.Range("K134").Value = 228000
$ws.Range("L134").Value = 172336.2
$ws.Range("M134").Value = -225465
$ws.Range("N134").Value = -177406.2
$ws.Range("H135").Value = 89290.25
$ws.Range("J135").Value = 89290.25
$ws.Range("L135").Value = 89290.25
$ws.Range("N135").Value = -99430.25
$ws.Range("H137").Value = 80876.44500000001
$ws.Range("J137").Value = 80876.44500000001
$ws.Range("L137").Value = 80876.44500000001
$ws.Range("N137").Value = -91076.44500000001
